# Applies the "Week 16 logged, season sim from Week 17" stat updates.
$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$rushing = $wb.Worksheets.Item("Rushing")

# J.Burrow (row 2): 3DATT 11 -> 12
$rushing.Range("E2").Value = 12

# J.Mixon (row 4): 1DATT 165 -> 178, 2DATT 81 -> 84, 3DATT 14 -> 16
$rushing.Range("C4").Value = 178
$rushing.Range("D4").Value = 84
$rushing.Range("E4").Value = 16

# C.Evans (row 6): 1DATT 5 -> 6
$rushing.Range("C6").Value = 6

# --- Receiving sheet ---
$receiving = $wb.Worksheets.Item("Receiving")

# J.Mixon (row 2): Short Target 32 -> 37, Short Comp 27 -> 32, Deep Target 1 -> 2, Deep Comp 1 -> 2
$receiving.Range("C2").Value = 37
$receiving.Range("D2").Value = 32
$receiving.Range("E2").Value = 2
$receiving.Range("F2").Value = 2

# S.Perine (row 3): Short Target 27 -> 29, Short Comp 23 -> 25
$receiving.Range("C3").Value = 29
$receiving.Range("D3").Value = 25

# C.Evans (row 4): Short Target 10 -> 11
$receiving.Range("C4").Value = 11

# T.Boyd (row 5): Short Target 67 -> 70, Short Comp 51 -> 53, Deep Target 17 -> 19, Deep Comp 10 -> 11
$receiving.Range("C5").Value = 70
$receiving.Range("D5").Value = 53
$receiving.Range("E5").Value = 19
$receiving.Range("F5").Value = 11

# T.Higgins (row 6): Short Target 77 -> 86, Short Comp 54 -> 63, Deep Target 27 -> 30, Deep Comp 16 -> 19
$receiving.Range("C6").Value = 86
$receiving.Range("D6").Value = 63
$receiving.Range("E6").Value = 30
$receiving.Range("F6").Value = 19

# J.Chase (row 7): Short Target 70 -> 77, Short Comp 47 -> 52, Deep Target 33 -> 36, Deep Comp 15 -> 17
$receiving.Range("C7").Value = 77
$receiving.Range("D7").Value = 52
$receiving.Range("E7").Value = 36
$receiving.Range("F7").Value = 17

# C.Uzomah (row 11): Short Target 42 -> 49, Short Comp 33 -> 38
$receiving.Range("C11").Value = 49
$receiving.Range("D11").Value = 38

# D.Sample (row 12): Short Target 11 -> 13, Short Comp 8 -> 10
$receiving.Range("C12").Value = 13
$receiving.Range("D12").Value = 10
